$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3890.5881
$ws.Range("I40").Value = 2332
$ws.Range("J40").Value = 4224.5713
$ws.Range("K40").Value = 2332
$ws.Range("L40").Value = 4224.5713
$ws.Range("M40").Value = -2157
$ws.Range("N40").Value = -4574.5713
$ws.Range("H101").Value = 524.7273
$ws.Range("I101").Value = 547
$ws.Range("K101").Value = 1641
$ws.Range("M101").Value = -19
$ws.Range("H107").Value = 1009.1667
$ws.Range("J107").Value = 941.6667
$ws.Range("L107").Value = 941.6667
$ws.Range("N107").Value = -4781.6667
$ws.Range("H111").Value = 9650.571
$ws.Range("I111").Value = 12045.625
$ws.Range("J111").Value = 4425
$ws.Range("K111").Value = 36136.875
$ws.Range("L111").Value = 13275
$ws.Range("M111").Value = -33069.875
$ws.Range("N111").Value = -19409
$ws.Range("H115").Value = 1122.7778
$ws.Range("I115").Value = 1122.7778
$ws.Range("K115").Value = 3368.3334
$ws.Range("M115").Value = -1801.3334
$ws.Range("H131").Value = 3338.5
$ws.Range("I131").Value = 2732.9092
$ws.Range("K131").Value = 8198.7276
$ws.Range("M131").Value = -3158.7276
$ws.Range("H133").Value = 41499.25
$ws.Range("J133").Value = 41499.25
$ws.Range("L133").Value = 41499.25
$ws.Range("N133").Value = -51619.25
$ws.Range("H136").Value = 42666.668
$ws.Range("J136").Value = 42666.668
$ws.Range("L136").Value = 42666.668
$ws.Range("N136").Value = -52866.668
$ws.Range("H138").Value = 187540.39
$ws.Range("I138").Value = 940.3125
$ws.Range("J138").Value = 266108.84
$ws.Range("K138").Value = 2820.9375
$ws.Range("L138").Value = 798326.52
$ws.Range("M138").Value = 2319.0625
$ws.Range("N138").Value = -808606.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7860.9453
$ws.Range("I32").Value = 7860.9453
$ws.Range("K32").Value = 7860.9453
$ws.Range("M32").Value = -7573.9453
$ws.Range("H61").Value = 11247.28
$ws.Range("I61").Value = 9514.076999999999
$ws.Range("K61").Value = 9514.076999999999
$ws.Range("M61").Value = -9302.076999999999
$ws.Range("H74").Value = 2363.7097
$ws.Range("I74").Value = 1143.375
$ws.Range("K74").Value = 1143.375
$ws.Range("M74").Value = -269.375
$ws.Range("H77").Value = 2363.7097
$ws.Range("I77").Value = 1143.375
$ws.Range("K77").Value = 5716.875
$ws.Range("M77").Value = -1348.875
$ws.Range("H110").Value = 1848.4286
$ws.Range("I110").Value = 1789.2727
$ws.Range("J110").Value = 2065.3333
$ws.Range("K110").Value = 1789.2727
$ws.Range("L110").Value = 2065.3333
$ws.Range("M110").Value = 255.7273
$ws.Range("N110").Value = -6155.3333
$ws.Range("H132").Value = 4583.4
$ws.Range("I132").Value = 4583.4
$ws.Range("K132").Value = 13750.2
$ws.Range("M132").Value = -11220.2
$ws.Range("H136").Value = 11247.28
$ws.Range("I136").Value = 9514.076999999999
$ws.Range("K136").Value = 28542.231
$ws.Range("M136").Value = -25992.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 62077.09
$ws.Range("J35").Value = 62077.09
$ws.Range("L35").Value = 62077.09
$ws.Range("N35").Value = -62697.09
$ws.Range("H82").Value = 30757.47
$ws.Range("I82").Value = 5541.1665
$ws.Range("K82").Value = 5541.1665
$ws.Range("M82").Value = -5158.1665
$ws.Range("H85").Value = 30757.47
$ws.Range("I85").Value = 5541.1665
$ws.Range("K85").Value = 5541.1665
$ws.Range("M85").Value = -4215.1665
$ws.Range("H94").Value = 1415.4615
$ws.Range("I94").Value = 999.1111
$ws.Range("K94").Value = 999.1111
$ws.Range("M94").Value = -548.1111
$ws.Range("H99").Value = 2684.739
$ws.Range("I99").Value = 1600.1875
$ws.Range("K99").Value = 1600.1875
$ws.Range("M99").Value = -102.1875
$ws.Range("H134").Value = 6230.5
$ws.Range("I134").Value = 3513.2778
$ws.Range("J134").Value = 8453.682000000001
$ws.Range("K134").Value = 10539.8334
$ws.Range("L134").Value = 25361.046
$ws.Range("M134").Value = -8004.8334
$ws.Range("N134").Value = -30431.046

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 498.05554
$ws.Range("I16").Value = 525.9167
$ws.Range("J16").Value = 442.33334
$ws.Range("K16").Value = 525.9167
$ws.Range("L16").Value = 442.33334
$ws.Range("M16").Value = -238.9167
$ws.Range("N16").Value = -1016.33334
$ws.Range("H99").Value = 3102.2
$ws.Range("I99").Value = 2499.6667
$ws.Range("K99").Value = 2499.6667
$ws.Range("M99").Value = -1001.6667
$ws.Range("H113").Value = 498.05554
$ws.Range("I113").Value = 525.9167
$ws.Range("J113").Value = 442.33334
$ws.Range("K113").Value = 525.9167
$ws.Range("L113").Value = 442.33334
$ws.Range("M113").Value = 1644.0833
$ws.Range("N113").Value = -4782.33334
$ws.Range("H122").Value = 3297.3845
$ws.Range("I122").Value = 2983.4736
$ws.Range("J122").Value = 4149.4287
$ws.Range("K122").Value = 8950.4208
$ws.Range("L122").Value = 12448.2861
$ws.Range("M122").Value = -6500.4208
$ws.Range("N122").Value = -17348.2861
$ws.Range("H126").Value = 3102.2
$ws.Range("I126").Value = 2499.6667
$ws.Range("K126").Value = 7499.000100000001
$ws.Range("M126").Value = -5029.000100000001
$ws.Range("H134").Value = 2678.439
$ws.Range("I134").Value = 1394.3334
$ws.Range("J134").Value = 6180.5454
$ws.Range("K134").Value = 4183.0002
$ws.Range("L134").Value = 18541.6362
$ws.Range("M134").Value = -1648.0002
$ws.Range("N134").Value = -23611.6362
$ws.Range("H137").Value = 91478
$ws.Range("J137").Value = 91478
$ws.Range("L137").Value = 91478
$ws.Range("N137").Value = -101678

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3001.5
$ws.Range("J48").Value = 500
$ws.Range("L48").Value = 1500
$ws.Range("N48").Value = -2000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1729.7273
$ws.Range("I113").Value = 1804.1111
$ws.Range("K113").Value = 1804.1111
$ws.Range("M113").Value = 365.8888999999999
$ws.Range("H122").Value = 3916.2942
$ws.Range("I122").Value = 2838.7
$ws.Range("K122").Value = 8516.099999999999
$ws.Range("M122").Value = -6066.099999999999
$ws.Range("H126").Value = 4034.375
$ws.Range("J126").Value = 3972
$ws.Range("L126").Value = 11916
$ws.Range("N126").Value = -16856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3734.1707
$ws.Range("I7").Value = 2704.0356
$ws.Range("J7").Value = 5952.923
$ws.Range("K7").Value = 2704.0356
$ws.Range("L7").Value = 5952.923
$ws.Range("M7").Value = -2592.0356
$ws.Range("N7").Value = -6176.923
$ws.Range("H40").Value = 4178.5386
$ws.Range("I40").Value = 3888.9565
$ws.Range("K40").Value = 3888.9565
$ws.Range("M40").Value = -3752.9565
$ws.Range("H46").Value = 10250
$ws.Range("I46").Value = 10725
$ws.Range("J46").Value = 10197.223
$ws.Range("K46").Value = 10725
$ws.Range("L46").Value = 10197.223
$ws.Range("M46").Value = -10537
$ws.Range("N46").Value = -10573.223
$ws.Range("H93").Value = 2132.8572
$ws.Range("I93").Value = 2146.4
$ws.Range("J93").Value = 2099
$ws.Range("K93").Value = 2146.4
$ws.Range("L93").Value = 2099
$ws.Range("M93").Value = -898.4000000000001
$ws.Range("N93").Value = -4595
$ws.Range("H126").Value = 3734.1707
$ws.Range("I126").Value = 2704.0356
$ws.Range("J126").Value = 5952.923
$ws.Range("K126").Value = 8112.1068
$ws.Range("L126").Value = 17858.769
$ws.Range("M126").Value = -5642.1068
$ws.Range("N126").Value = -22798.769
$ws.Range("H132").Value = 3065.3276
$ws.Range("I132").Value = 3047.1738
$ws.Range("K132").Value = 9141.5214
$ws.Range("M132").Value = -6611.5214
$ws.Range("H136").Value = 4759.278
$ws.Range("I136").Value = 4371.8
$ws.Range("K136").Value = 13115.4
$ws.Range("M136").Value = -10565.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2612.1
$ws.Range("I122").Value = 2395.5334
$ws.Range("J122").Value = 3261.8
$ws.Range("K122").Value = 7186.600199999999
$ws.Range("L122").Value = 9785.400000000001
$ws.Range("M122").Value = -4736.600199999999
$ws.Range("N122").Value = -14685.4
$ws.Range("H126").Value = 2653.3157
$ws.Range("I126").Value = 2213.375
$ws.Range("K126").Value = 6640.125
$ws.Range("M126").Value = -4170.125
$ws.Range("H132").Value = 2250.4092
$ws.Range("I132").Value = 2194.9443
$ws.Range("K132").Value = 6584.8329
$ws.Range("M132").Value = -4054.8329
$ws.Range("H136").Value = 8797.814
$ws.Range("I136").Value = 9407.24
$ws.Range("J136").Value = 1180
$ws.Range("K136").Value = 28221.72
$ws.Range("L136").Value = 3540
$ws.Range("M136").Value = -25671.72
$ws.Range("N136").Value = -8640

